$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.466.49"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.569.26"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'208.80"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").Value = "'0.500"
$ws.Range("E6").Value = "  -1.21%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D10").Value = "'0.0593"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "1.793.60"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").Value = "1.579.36"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("E15").Value = "  -2.46%  "
$ws.Range("D16").Value = "'63.77"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "27.481.73"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "'214.34"
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").Value = "'7.29"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "'4.12"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "'9.56"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "'2.01"
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("D25").Value = "'152.85"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").Value = "'15.02"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("D32").Value = "'3.20"
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").Value = "1.380.87"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("E34").Value = "  +1.86%  "
$ws.Range("D35").Value = "'1.54"
$ws.Range("E35").Value = "  +1.09%  "
$ws.Range("D36").Value = "'2.31"
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("D37").Value = "'0.953"
$ws.Range("E37").Value = "  -2.61%  "
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("D39").Value = "'0.544"
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("D40").Value = "'0.827"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").Value = "'0.983"
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("D43").Value = "'1.81"
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("D44").Value = "'64.25"
$ws.Range("E44").Value = "  +1.04%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("D47").Value = "1.705.25"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").Value = "'85.37"
$ws.Range("E48").Value = "  -3.40%  "
$ws.Range("D49").Value = "0.0₆0100"
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").Value = "'0.0959"
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("E51").Value = "  -0.56%  "

$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D50").ClearFormats()
